# Updates the cryptos price/volume table (columns D and E) with latest
# scraped values, as produced by the scheduled "Updated cryptos list" GitHub Action.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.411.18'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.01%  '

# Row 3: Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.850.90'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.20%  '

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.77'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.03%  '

# Row 6: XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6290'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.35%  '

# Row 7: USDC
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.02%  '

# Row 8: Dogecoin
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07677'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.65%  '

# Row 9: Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2934'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.62%  '

# Row 10: Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.55'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.21%  '

# Row 11: TRON
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.70%  '

# Row 12: WrappedEther
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.853.51'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.17%  '

# Row 13: ShibaInu
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.00001107'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +10.20%  '

# Row 14: Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.024'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.74%  '

# Row 16: Litecoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.62'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.71%  '

# Row 17: WrappedliquidstakedEther2.0
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.105.39'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.10%  '

# Row 18: Uniswap
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.149'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.44%  '

# Row 19: WrappedBTC
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '29.453.84'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.07%  '

# Row 20: BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '229.03'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.23%  '

# Row 21: Avalanche
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.46'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.10%  '

# Row 22: Dai
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.05%  '

# Row 23: Chainlink
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.432'

# Row 24: BinanceUSD
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.001'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.01%  '

# Row 25: Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.03'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.02%  '

# Row 26: Stellar
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1385'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.57%  '

# Row 27: Cosmos
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.389'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.22%  '

# Row 28: EthereumClassic
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.16%  '

# Row 29: Toncoin
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +4.11%  '

# Row 30: PancakeSwap
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.467'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.12%  '

# Row 31: Hedera
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.05712'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.25%  '

# Row 32: Filecoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.127'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.14%  '

# Row 33: InternetComputer(DFINITY)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.050'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.67%  '

# Row 34: LidoDAOToken
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.15%  '

# Row 35: ARBITRUM
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.164'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.79%  '

# Row 36: ImmutableX
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7078'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.95%  '

# Row 37: HuobiToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.585'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.10%  '

# Row 38: MXToken
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.00%  '

# Row 39: VeChain
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01790'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.03%  '

# Row 40: Maker
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.218.94'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.34%  '

# Row 41: FraxShare
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.481'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +5.01%  '

# Row 42: TrustWalletToken
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9100'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.08%  '

# Row 43: PaxDollar
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.02%  '

# Row 44: RocketPoolETH
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.014.12'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.11%  '

# Row 45: Quant
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.71'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.07%  '

# Row 46: Aave
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '66.30'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.41%  '

# Row 47: Aptos
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.123'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.48%  '

# Row 48: BabyDogeCoin
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000118'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.75%  '

# Row 49: TheSandbox
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4015'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.16%  '

# Row 50: EnergySwap
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.979'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.11%  '

# Row 51: RenderToken
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.11%  '
